# Actualización automática 2025-08-27 15:05:09
# Insert a new client row ("FRANK FERRETERIA FRANKFERRE CIA.") right above
# the existing "VIEJO RIVAS MAYRA ANABELLE" row on both sheets, pushing the
# summary/totals row down and updating its "0 de N" label on the first sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" (columns A:R) ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(6).Insert()

$ws1.Cells.Item(6, 1).Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws1.Cells.Item(6, 2).Value = "FRANK FERRETERIA FRANKFERRE CIA."
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(6, $col).Value = 0
}

# Update the totals row (now row 8) label from "0 de 5" to "0 de 6"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(8, $col).Value = "0 de 6"
}

# --- Sheet "VENTA MENSUAL" (columns A:G) ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(6).Insert()

$ws2.Cells.Item(6, 1).Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws2.Cells.Item(6, 2).Value = "FRANK FERRETERIA FRANKFERRE CIA."
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(6, $col).Value = 0
}
